# Weekly fruit/vegetable price refresh: two new daily observations were
# appended to the "Berenjena" (Terminal Hortofrutícola Agro Chillán) series.
# Appending new rows to this date-ordered table is implemented as inserting
# rows, which pushes the existing rows below the insertion point down by one:
#
#   - New record for 2023-07-13 is inserted at row 72, pushing the old
#     rows 72-97 down to rows 73-98.
#   - New record for 2023-07-14 is inserted at row 98 (i.e. right before the
#     old row 97, which by then sits at row 98), pushing it down to row 99.
#
# Net effect: dimension grows from A1:R97 to A1:R99.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row at 72 (Fecha serial 45120 = 2023-07-13) ---
$ws.Rows.Item(72).Insert()

$ws.Cells.Item(72, 1).Value = 7
$ws.Cells.Item(72, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(72, 3).Value = "Ñuble"
$ws.Cells.Item(72, 4).Value = 45120
$ws.Cells.Item(72, 5).Value = 16
$ws.Cells.Item(72, 6).Value = 100112001
$ws.Cells.Item(72, 7).Value = "Berenjena"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 80
$ws.Cells.Item(72, 11).Value = 8000
$ws.Cells.Item(72, 12).Value = 8000
$ws.Cells.Item(72, 13).Value = 8000
$ws.Cells.Item(72, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(72, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(72, 16).Value = 133
$ws.Cells.Item(72, 17).Value = 60
$ws.Cells.Item(72, 18).Value = "Hortaliza"

# --- Insert new row at 98 (Fecha serial 45121 = 2023-07-14) ---
$ws.Rows.Item(98).Insert()

$ws.Cells.Item(98, 1).Value = 7
$ws.Cells.Item(98, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(98, 3).Value = "Ñuble"
$ws.Cells.Item(98, 4).Value = 45121
$ws.Cells.Item(98, 5).Value = 16
$ws.Cells.Item(98, 6).Value = 100112001
$ws.Cells.Item(98, 7).Value = "Berenjena"
$ws.Cells.Item(98, 8).Value = "Sin especificar"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 80
$ws.Cells.Item(98, 11).Value = 7000
$ws.Cells.Item(98, 12).Value = 8000
$ws.Cells.Item(98, 13).Value = 7625
$ws.Cells.Item(98, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(98, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(98, 16).Value = 127
$ws.Cells.Item(98, 17).Value = 60
$ws.Cells.Item(98, 18).Value = "Hortaliza"
